$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: Quantity 70 -> 52
$ws.Range("D3").Value = 52

# Row 4: Quantity 0 -> 7, Date 2024-09-07 -> 2024-09-09, Time 17:33:16 -> 21:40:23
$ws.Range("D4").Value = 7
$ws.Range("F4").Value = "'2024-09-09"
$ws.Range("F4").Style = "Normal"
$ws.Range("G4").Value = "'21:40:23"
$ws.Range("G4").Style = "Normal"

# Row 6: Quantity 7 -> 6
$ws.Range("D6").Value = 6
